$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.656.17"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "1.530.44"
$ws.Range("E3").Value = "  -1.95%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'205.21"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").Value = "'0.483"
$ws.Range("E6").Value = "  -1.24%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -1.37%  "

$ws.Range("D9").Value = "'21.25"
$ws.Range("E9").Value = "  -3.13%  "

$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("E11").Value = "  -1.31%  "

$ws.Range("D12").Value = "1.748.15"
$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").Value = "1.540.93"
$ws.Range("E13").Value = "  -1.31%  "

$ws.Range("D14").Value = "'3.65"
$ws.Range("E14").Value = "  -2.16%  "

$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'61.22"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "26.661.23"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").Value = "'211.44"
$ws.Range("E18").Value = "  -1.57%  "

$ws.Range("D19").Value = "0.0₃0683"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("E20").Value = "  -2.55%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("E23").Value = "  -3.25%  "

$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("D25").Value = "'152.33"
$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("E26").Value = "  -3.66%  "

$ws.Range("D27").Value = "'14.82"
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("E30").Value = "  -1.82%  "

$ws.Range("E31").Value = "  -1.92%  "

$ws.Range("E32").Value = "  +2.23%  "

$ws.Range("D33").Value = "1.353.52"
$ws.Range("E33").Value = "  -2.42%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  -3.75%  "

$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("D37").Value = "'0.933"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("D39").Value = "'0.521"
$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("D40").Value = "'0.795"
$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("E41").Value = "  +4.81%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "'62.25"
$ws.Range("E44").Value = "  -1.62%  "

$ws.Range("E45").Value = "  -1.88%  "

$ws.Range("E46").Value = "  -4.01%  "

$ws.Range("D47").Value = "1.663.46"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("D48").Value = "'85.64"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("E49").Value = "  +3.23%  "

$ws.Range("D50").Value = "'0.0945"
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("E51").Value = "  +0.04%  "
